$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2 through 20
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
